$d = $word.ActiveDocument

# The paragraph "Preparing food with heat or fire is an activity unique to
# humans." sits in its own paragraph right after the paragraph ending in
# "...without the presence of heat,". Removing it should both delete that
# sentence and merge its paragraph back into the previous one (the trailing
# ". " run then becomes the tail of the previous paragraph).
$sentence = "Preparing food with heat or fire is an activity unique to humans."

$searchRange = $d.Content
$found = $searchRange.Find.Execute($sentence, $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the target sentence to remove."
}

$start = $searchRange.Start
$end = $searchRange.End

# 1) Delete the sentence text itself. This leaves the (now empty) paragraph
#    mark that used to separate it from the previous paragraph, followed by
#    the trailing ". " run.
$d.Range($start, $end).Delete()

# 2) Delete that paragraph mark so the remaining paragraph (just ". ") merges
#    back up into the previous paragraph ("...without the presence of heat,").
$d.Range($start - 1, $start).Delete()
